$d = $word.ActiveDocument

# 1. Merge the two runs "CD " and "C:\Users\andrew.dilley\development\sql-demo"
#    into a single run with the combined text (same Courier New formatting).
$d.Content.Find.Execute(
    "CD C:\Users\andrew.dilley\development\sql-demo",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "CD C:\Users\andrew.dilley\development\sql-demo", 2) | Out-Null

# 2. After the last paragraph ("git push -u origin main"), append five new
#    paragraphs (all using the same Courier New formatting as the
#    surrounding text): blank, blank, "March 7", blank, blank.
#    Note: when InsertAfter-ing at a Range collapsed to the very end of the
#    document, a single *leading* paragraph mark in the inserted text is
#    absorbed into the existing end-of-story mark rather than creating an
#    extra empty paragraph, so three leading marks are needed to get two
#    blank paragraphs before "March 7".
$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertAfter("`r`r`rMarch 7`r`r")
